$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")
Write-Output $ws.UsedRange.Address()
